$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New header labels for columns G (layerType) and H (dataType), plus a
# duplicate "dataType" header in the new column N.
# Order matters for shared-string indices: dataType must be interned
# before layerType so dataType gets the lower index (172) and layerType
# gets 173, matching the target workbook.
$ws.Range("H1").Value = "dataType"
$ws.Range("G1").Value = "layerType"
$ws.Range("N1").Value = "dataType"

# New column N holds a count value (derived from the dataType) for every
# data row already populated in the sheet.
$nValues = @{
    2 = 1
    3 = 1
    4 = 2
    5 = 2
    6 = 1
    7 = 1
    8 = 2
    9 = 3
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 2
    15 = 3
    16 = 3
    17 = 2
    18 = 1
    19 = 3
    20 = 2
    21 = 2
    22 = 2
    23 = 2
    24 = 2
    25 = 1
    26 = 2
    27 = 2
    28 = 2
    29 = 2
    30 = 2
    31 = 2
    32 = 2
    33 = 2
    34 = 2
    35 = 3
}

foreach ($row in $nValues.Keys) {
    $ws.Cells.Item($row, 14).Value = $nValues[$row]
}

# Column M (symbologyType's long source text) is now hidden.
$ws.Columns("M:M").Hidden = $true

# Move the active selection / reset the scroll position.
$ws.Range("H13").Select()
